$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- 1. Cell shading: add explicit w:color="auto" w:fill="auto" to every
#        score-table cell that currently has a bare <w:shd w:val="clear"/>.
#        Rows 6..11 (name cell in col 2, score cell in col 4) all get it.
for ($r = 6; $r -le 11; $r++) {
    $nameShd = $tbl.Cell($r, 2).Shading
    $nameShd.ForegroundPatternColor = -16777216   # wdColorAutomatic -> w:color="auto"
    $nameShd.BackgroundPatternColor = -16777216   # wdColorAutomatic -> w:fill="auto"

    $scoreShd = $tbl.Cell($r, 4).Shading
    $scoreShd.ForegroundPatternColor = -16777216
    $scoreShd.BackgroundPatternColor = -16777216
}

# --- 2. Update the score values. Assigning Range.Text directly (instead of
#        Find/Replace) keeps each edit scoped to its own cell and preserves
#        the run's existing formatting.
$tbl.Cell(6, 4).Range.Text = "8.4"    # 白靖妍: 8.5 -> 8.4
$tbl.Cell(7, 4).Range.Text = "9.5"    # 李思涵: 8.5 -> 9.5
# row 8 (马雯丽: 9.0) is unchanged
$tbl.Cell(9, 4).Range.Text = "9.6"    # 赵益萍: 9.1 -> 9.6
$tbl.Cell(10, 4).Range.Text = "9.6"   # 仵梦雅: 9.0 -> 9.6
$tbl.Cell(11, 4).Range.Text = "8.5"   # 王佳丽: 9.5 -> 8.5

# --- 3. Move the hidden "_GoBack" bookmark from 仵梦雅's score paragraph
#        (row 10) to 白靖妍's score paragraph (row 6), right after the text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$tbl2 = $d.Tables.Item(1)
$targetRange = $tbl2.Cell(6, 4).Range
$bmPos = $targetRange.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
